$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.094.28'
$ws.Range('D3').Value = '1.650.83'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '''218.29'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').Value = '''0.5204'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = '''1.003'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '''0.2642'
$ws.Range('E8').Value = '  +0.80%  '
$ws.Range('D9').Value = '''0.06330'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').Value = '''20.38'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = '''0.07694'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '''4.601'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('D13').Value = '1.654.34'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').Value = '1.878.47'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '''0.5597'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '0.0₅8150'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').Value = '''65.35'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '26.094.45'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '''4.622'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = '''10.45'
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('D22').Value = '''191.27'
$ws.Range('E22').Value = '  -1.65%  '
$ws.Range('D23').Value = '''5.919'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '''144.21'
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '''7.224'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').Value = '''15.91'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '''0.05485'
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('D31').Value = '''1.267'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = '''3.443'
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').Value = '''3.365'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').Value = '''1.558'
$ws.Range('E34').Value = '  -2.26%  '
$ws.Range('D35').Value = '''0.9486'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').Value = '''2.779'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = '''2.401'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').Value = '''0.5638'
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('D39').Value = '''0.01578'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').Value = '''5.858'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.032.74'
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''0.8311'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Value = '''101.19'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('D45').Value = '1.791.49'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''57.54'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈108'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('D48').Value = '''0.9995'
$ws.Range('D49').Value = '''0.4338'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = '''7.966'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '''0.05165'
$ws.Range('E51').Value = '  -2.42%  '
